# Commit: "commiting the test cases"
# Update the "To Be Executed" (column C) answers for the Testcases sheet from
# "Y" to "N" for rows 19-77 (row 18 stays "Y"), and move the active
# selection to B74 to match where the author ended up after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

# Flip column C (rows 19 through 77) from "Y" to "N".
$ws.Range("C19:C77").Value2 = "N"

# Match the final selection left behind in the saved file.
$ws.Range("B74").Select()
